# Update generated output numbers (gh-pages data refresh at 456a3b4)
$wb = $excel.ActiveWorkbook

# Sheet "展览" - column F values bump up slightly
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 5536
$wsExhibition.Range("F13").Value = 5082
$wsExhibition.Range("F20").Value = 4363
$wsExhibition.Range("F23").Value = 120

# Sheet "全部类型" - same underlying rows, offset by one row
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 5536
$wsAll.Range("F14").Value = 5082
$wsAll.Range("F21").Value = 4363
$wsAll.Range("F24").Value = 120
